# "testing CS with look back 1 year"
# Update a handful of cells on the "Coupling Parameters" sheet and move the
# active selection / view to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# End Year: 2055 -> 2089
$ws.Range("B3").Value = 2089

# change_IRR: TRUE -> FALSE
$ws.Range("B5").Value = $false

# dynamic_derating_factor: TRUE -> FALSE
$ws.Range("B7").Value = $false

# scenarioWeatheryearsExcel: 40weatherYears2050TNO.xlsx -> 40weatherYears2050TNO-2004.xlsx
$ws.Range("B31").Value = "40weatherYears2050TNO-2004.xlsx"

# Scroll the view down so row 42 is at the top, and move the selection to B7,
# matching where the author ended up when they saved the workbook.
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
